$wb = $excel.ActiveWorkbook

# Sheet "o": updated simulation output values
$ws = $wb.Worksheets.Item("o")
$ws.Range("A1").Value = [double]"0.7632407081291448"
$ws.Range("B1").Value = [double]"0.01461892953374144"
$ws.Range("A2").Value = [double]"0.1776590689724449"
$ws.Range("B2").Value = [double]"0.728158415035312"
$ws.Range("C2").Value = [double]"0.6955074434115286"
$ws.Range("D2").Value = [double]"0.1849920720800949"
$ws.Range("A3").Value = [double]"8.112697552874444E-18"
$ws.Range("B3").Value = [double]"3.317029860677889E-17"
$ws.Range("C3").Value = [double]"3.141128805256694E-17"
$ws.Range("D3").Value = [double]"8.00434827190673E-18"
$ws.Range("A4").Value = [double]"5.729642221079343E-17"
$ws.Range("B4").Value = [double]"2.342726833367584E-16"
$ws.Range("C4").Value = [double]"2.218674835237584E-16"
$ws.Range("D4").Value = [double]"5.655993401404293E-17"
$ws.Range("A5").Value = [double]"0.05910022289841008"
$ws.Range("B5").Value = [double]"0.2572226554309462"
$ws.Range("C5").Value = [double]"0.3044925565884709"
$ws.Range("D5").Value = [double]"0.815007927919905"

# Sheet "z": canola stubble cat A size change to 1%, grain_propn rework
$ws = $wb.Worksheets.Item("z")
$ws.Range("A1").Value = [double]"0.7168348031366673"
$ws.Range("B1").Value = [double]"0.02978849117830595"
$ws.Range("A2").Value = [double]"0.2823222605127404"
$ws.Range("B2").Value = [double]"0.9665379169450865"
$ws.Range("C2").Value = [double]"0.3888902535642045"
$ws.Range("A3").Value = [double]"1.090492102067622E-06"
$ws.Range("B3").Value = [double]"4.752492256304215E-06"
$ws.Range("C3").Value = [double]"0.0007935567444627874"
$ws.Range("D3").Value = [double]"0.001410268770281795"
$ws.Range("A4").Value = [double]"0.0001456984526927458"
$ws.Range("B4").Value = [double]"0.0006349474845136278"
$ws.Range("C4").Value = [double]"0.1037261208925042"
$ws.Range("D4").Value = [double]"0.0983035723041345"
$ws.Range("A5").Value = [double]"0.000696147405797672"
$ws.Range("B5").Value = [double]"0.003033891899837479"
$ws.Range("C5").Value = [double]"0.5065900687988284"
$ws.Range("D5").Value = [double]"0.9002861589255836"

# Sheet "r": same recalculated values as sheet "z"
$ws = $wb.Worksheets.Item("r")
$ws.Range("A1").Value = [double]"0.7168348031366673"
$ws.Range("B1").Value = [double]"0.02978849117830595"
$ws.Range("A2").Value = [double]"0.2823222605127404"
$ws.Range("B2").Value = [double]"0.9665379169450865"
$ws.Range("C2").Value = [double]"0.3888902535642045"
$ws.Range("A3").Value = [double]"1.090492102067622E-06"
$ws.Range("B3").Value = [double]"4.752492256304215E-06"
$ws.Range("C3").Value = [double]"0.0007935567444627874"
$ws.Range("D3").Value = [double]"0.001410268770281795"
$ws.Range("A4").Value = [double]"0.0001456984526927458"
$ws.Range("B4").Value = [double]"0.0006349474845136278"
$ws.Range("C4").Value = [double]"0.1037261208925042"
$ws.Range("D4").Value = [double]"0.0983035723041345"
$ws.Range("A5").Value = [double]"0.000696147405797672"
$ws.Range("B5").Value = [double]"0.003033891899837479"
$ws.Range("C5").Value = [double]"0.5065900687988284"
$ws.Range("D5").Value = [double]"0.9002861589255836"
